$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2-10 for columns I (I0) and J (IF)
$values = @{
    2 = @(8, 8)
    3 = @(8, 8)
    4 = @(9, 9)
    5 = @(9, 9)
    6 = @(9, 9)
    7 = @(8, 8)
    8 = @(7, 7)
    9 = @(4, 5)
    10 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
